$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 11968
$ws.Range("I46").Value = 17716.285
$ws.Range("J46").Value = 6219.7144
$ws.Range("K46").Value = 53148.855
$ws.Range("L46").Value = 18659.1432
$ws.Range("M46").Value = -53029.855
$ws.Range("N46").Value = -18897.1432

$ws.Range("H60").Value = 11968
$ws.Range("I60").Value = 17716.285
$ws.Range("J60").Value = 6219.7144
$ws.Range("K60").Value = 53148.855
$ws.Range("L60").Value = 18659.1432
$ws.Range("M60").Value = -52664.855
$ws.Range("N60").Value = -19627.1432

$ws.Range("H126").Value = 46951.2
$ws.Range("J126").Value = 46951.2
$ws.Range("L126").Value = 46951.2
$ws.Range("N126").Value = -56831.2

$ws.Range("H133").Value = 23274.643
$ws.Range("J133").Value = 23274.643
$ws.Range("L133").Value = 23274.643
$ws.Range("N133").Value = -33394.643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 5252.5
$ws.Range("I15").Value = 3011
$ws.Range("J15").Value = 5999.6665
$ws.Range("K15").Value = 3011
$ws.Range("L15").Value = 5999.6665
$ws.Range("M15").Value = -2661
$ws.Range("N15").Value = -6699.6665

$ws.Range("H34").Value = 35000
$ws.Range("J34").Value = 35000
$ws.Range("L34").Value = 35000
$ws.Range("N34").Value = -35542

$ws.Range("H80").Value = 49243.855
$ws.Range("J80").Value = 49243.855
$ws.Range("L80").Value = 49243.855
$ws.Range("N80").Value = -51239.855

$ws.Range("H83").Value = 49243.855
$ws.Range("J83").Value = 49243.855
$ws.Range("L83").Value = 147731.565
$ws.Range("N83").Value = -157715.565

$ws.Range("H123").Value = 35614.5
$ws.Range("J123").Value = 35614.5
$ws.Range("L123").Value = 35614.5
$ws.Range("N123").Value = -45414.5

$ws.Range("H125").Value = 48930.668
$ws.Range("J125").Value = 48930.668
$ws.Range("L125").Value = 48930.668
$ws.Range("N125").Value = -58770.668

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

$ws.Range("H130").Value = 47429
$ws.Range("J130").Value = 47429
$ws.Range("L130").Value = 47429
$ws.Range("N130").Value = -57469

$ws.Range("H131").Value = 49607
$ws.Range("J131").Value = 49607
$ws.Range("L131").Value = 49607
$ws.Range("N131").Value = -59687

$ws.Range("H138").Value = 41650
$ws.Range("J138").Value = 41650
$ws.Range("L138").Value = 41650
$ws.Range("N138").Value = -51930

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 49996
$ws.Range("J124").Value = 49996
$ws.Range("L124").Value = 49996
$ws.Range("N124").Value = -59816

$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620

$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

$ws.Range("H130").Value = 46199
$ws.Range("J130").Value = 46199
$ws.Range("L130").Value = 46199
$ws.Range("N130").Value = -56239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1340

$ws.Range("H20").Value = 44790.668
$ws.Range("J20").Value = 44790.668
$ws.Range("L20").Value = 44790.668
$ws.Range("N20").Value = -45262.668

$ws.Range("H30").Value = 44790.668
$ws.Range("J30").Value = 44790.668
$ws.Range("L30").Value = 44790.668
$ws.Range("N30").Value = -44972.668

$ws.Range("H128").Value = 44790.668
$ws.Range("J128").Value = 44790.668
$ws.Range("L128").Value = 44790.668
$ws.Range("N128").Value = -54750.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 5249.6665
$ws.Range("J100").Value = 5249.6665
$ws.Range("L100").Value = 15748.9995
$ws.Range("N100").Value = -17370.9995

$ws.Range("H131").Value = 2887.3872
$ws.Range("I131").Value = 10374.3
$ws.Range("J131").Value = 1447.5962
$ws.Range("K131").Value = 31122.9
$ws.Range("L131").Value = 4342.7886
$ws.Range("M131").Value = -26082.9
$ws.Range("N131").Value = -14422.7886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1641.7368
$ws.Range("I113").Value = 1658
$ws.Range("J113").Value = 1623.6666
$ws.Range("K113").Value = 1658
$ws.Range("L113").Value = 1623.6666
$ws.Range("M113").Value = 512
$ws.Range("N113").Value = -5963.6666

$ws.Range("H119").Value = 47761
$ws.Range("J119").Value = 47761
$ws.Range("L119").Value = 47761
$ws.Range("N119").Value = -57437

$ws.Range("H130").Value = 45784
$ws.Range("J130").Value = 45784
$ws.Range("L130").Value = 45784
$ws.Range("N130").Value = -55824

$ws.Range("H138").Value = 40963.637
$ws.Range("J138").Value = 40963.637
$ws.Range("L138").Value = 40963.637
$ws.Range("N138").Value = -51243.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1218
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1920
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1920
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -2510

$ws.Range("H27").Value = 1218
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1920
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1920
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -2134

$ws.Range("H93").Value = 3101.375
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 2830.1428
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 2830.1428
$ws.Range("M93").Value = -3752
$ws.Range("N93").Value = -5326.1428

$ws.Range("H127").Value = 50617.75
$ws.Range("J127").Value = 50617.75
$ws.Range("L127").Value = 50617.75
$ws.Range("N127").Value = -60537.75

$ws.Range("H130").Value = 37996
$ws.Range("J130").Value = 37996
$ws.Range("L130").Value = 37996
$ws.Range("N130").Value = -48036

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H93").Value = 36998.668
$ws.Range("J93").Value = 36998.668
$ws.Range("L93").Value = 36998.668
$ws.Range("N93").Value = -41990.668

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""

$ws.Range("H120").Value = 35206
$ws.Range("J120").Value = 35206
$ws.Range("L120").Value = 35206
$ws.Range("N120").Value = -44882

$ws.Range("H128").Value = 49715
$ws.Range("J128").Value = 49715
$ws.Range("L128").Value = 49715
$ws.Range("N128").Value = -59675

$ws.Range("H133").Value = 86819
$ws.Range("J133").Value = 86819
$ws.Range("L133").Value = 86819
$ws.Range("N133").Value = -96939

$ws.Range("H135").Value = 22636.117
$ws.Range("J135").Value = 22636.117
$ws.Range("L135").Value = 22636.117
$ws.Range("N135").Value = -32776.117

$ws.Range("H137").Value = 46274.918
$ws.Range("J137").Value = 46274.918
$ws.Range("L137").Value = 46274.918
$ws.Range("N137").Value = -56474.918
